$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.845.87"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "1.833.64"
$ws.Range("E3").Value = "  +1.02%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "230.85"
$ws.Range("E5").Value = "  -1.08%  "

$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("D8").Value = "39.52"
$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("D9").Value = "0.323"
$ws.Range("E9").Value = "  +1.21%  "

$ws.Range("D10").Value = "0.0683"
$ws.Range("E10").Value = "  -0.44%  "

$ws.Range("E11").Value = "  -1.23%  "

$ws.Range("D12").Value = "2.094.36"
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "11.32"
$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.817.92"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("D15").Value = "0.667"
$ws.Range("E15").Value = "  +1.24%  "

$ws.Range("D16").Value = "4.65"
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("D17").Value = "34.821.70"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").Value = "69.42"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").Value = "0.0₃0786"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("D20").Value = "239.55"
$ws.Range("E20").Value = "  +0.60%  "

$ws.Range("D21").Value = "12.20"
$ws.Range("E21").Value = "  +3.06%  "

$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").Value = "2.25"
$ws.Range("E24").Value = "  -1.26%  "

$ws.Range("D25").Value = "172.37"
$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("D26").Value = "7.74"
$ws.Range("E26").Value = "  -1.27%  "

$ws.Range("D27").Value = "0.124"
$ws.Range("E27").Value = "  +2.64%  "

$ws.Range("D28").Value = "17.34"
$ws.Range("E28").Value = "  -0.57%  "

$ws.Range("D29").Value = "1.51"
$ws.Range("E29").Value = "  -8.21%  "

$ws.Range("E30").Value = "  +0.21%  "

$ws.Range("D31").Value = "0.0549"
$ws.Range("E31").Value = "  -0.88%  "

$ws.Range("D32").Value = "3.90"
$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("D33").Value = "3.93"
$ws.Range("E33").Value = "  -1.27%  "

$ws.Range("E34").Value = "  +8.36%  "

$ws.Range("E35").Value = "  +3.10%  "

$ws.Range("E36").Value = "  +10.58%  "

$ws.Range("D37").Value = "0.699"
$ws.Range("E37").Value = "  +2.85%  "

$ws.Range("D38").Value = "91.68"
$ws.Range("E38").Value = "  -1.56%  "

$ws.Range("D39").Value = "1.342.07"

$ws.Range("D40").Value = "1.03"
$ws.Range("E40").Value = "  +4.12%  "

$ws.Range("D41").Value = "0.0194"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").Value = "14.53"
$ws.Range("E42").Value = "  -1.55%  "

$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("E44").Value = "  -3.30%  "

$ws.Range("D45").Value = "2.76"
$ws.Range("E45").Value = "  -0.27%  "

$ws.Range("D46").Value = "6.26"
$ws.Range("E46").Value = "  -0.58%  "

$ws.Range("D47").Value = "0.0522"
$ws.Range("E47").Value = "  +1.88%  "

$ws.Range("D48").Value = "2.010.50"
$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").Value = "0.0671"
$ws.Range("E50").Value = "  +4.30%  "

$ws.Range("D51").Value = "3.22"
$ws.Range("E51").Value = "  +13.53%  "

